$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.193.15"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.736.49"
$ws.Range("E3").Value = "  +7.01%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "'519.52"
$ws.Range("E5").Value = "  +2.92%  "
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "2.732.14"
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").Value = "'6.30"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  +7.52%  "
$ws.Range("E12").Value = "  +3.07%  "
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "3.165.15"
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("D15").Value = "59.157.09"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("D16").Value = "'21.31"
$ws.Range("E16").Value = "  +4.13%  "
$ws.Range("E17").Value = "  +3.79%  "
$ws.Range("D18").Value = "2.722.80"
$ws.Range("E18").Value = "  +6.38%  "
$ws.Range("D19").Value = "'357.00"
$ws.Range("E19").Value = "  +7.63%  "
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E21").Value = "  +5.39%  "
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'61.45"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E25").Value = "  +5.63%  "
$ws.Range("D26").Value = "2.804.98"
$ws.Range("E26").Value = "  +4.95%  "
$ws.Range("D27").Value = "'0.163"
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("D28").Value = "'0.989"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "0.0₃0827"
$ws.Range("E29").Value = "  +6.58%  "
$ws.Range("E30").Value = "  +6.62%  "
$ws.Range("D31").Value = "'0.996"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "'6.46"
$ws.Range("E32").Value = "  +12.02%  "
$ws.Range("D33").Value = "'19.23"
$ws.Range("E33").Value = "  +4.04%  "
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").Value = "'150.38"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  +15.98%  "
$ws.Range("E37").Value = "  +6.45%  "
$ws.Range("E38").Value = "  +5.80%  "
$ws.Range("D39").Value = "'0.866"
$ws.Range("E39").Value = "  +5.92%  "
$ws.Range("D40").Value = "'37.00"
$ws.Range("E40").Value = "  +3.40%  "
$ws.Range("E41").Value = "  +7.04%  "
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").Value = "'0.630"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "'283.67"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'20.42"
$ws.Range("E45").Value = "  +10.08%  "
$ws.Range("D46").Value = "'0.0989"
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").Value = "'0.992"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("D49").Value = "'4.81"
$ws.Range("E49").Value = "  +6.89%  "
$ws.Range("E50").Value = "  +3.38%  "
$ws.Range("D51").Value = "2.025.53"
